# Daily attendance processing - 2026-01-18 15:58:31
#
# In the "Recorded By" column (G) of the "Session Analysis Results" sheet,
# swap the order of the two recorder names so that
#   "dnasr281@gmail.com, System"
# becomes
#   "System, dnasr281@gmail.com"
# for every row where that exact text currently appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
